# Update "Lũy kế tháng LONG XUYÊN" report rows with refreshed Notion data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newEditedTime = "2024-07-24T16:01:00.000Z"

$updates = @(
    @{ Row = 3;  T = 23500000; AA = 61142000; AE = 121700000; AH = 104000000 },
    @{ Row = 4;  T = 27500000; AA = 55199000; AE = 135400000; AH = 118400000 },
    @{ Row = 5;  T = 5500000;  AA = 17626000; AE = 45950000;  AH = 42650000  },
    @{ Row = 7;  T = 54500000; AA = 29992000; AE = 91000000;  AH = 91000000  },
    @{ Row = 13; T = 24000000; AA = 71500000; AE = 74000000;  AH = 67000000  }
)

foreach ($u in $updates) {
    $row = $u.Row

    $ws.Range("D$row").Value = $newEditedTime
    $ws.Range("T$row").Value = $u.T
    $ws.Range("AA$row").Value = $u.AA
    $ws.Range("AE$row").Value = $u.AE
    $ws.Range("AH$row").Value = $u.AH
}
